$wb = $excel.ActiveWorkbook

# --- Sheet "Mes actual cel" (sheet1) ---
$wsCel = $wb.Worksheets.Item("Mes actual cel")
$wsCel.Range("B2").Value = "Jesús Gutiérrez Vázquez"
$wsCel.Range("B3").Value = "Jose Luis Pérez Asenjo"
$wsCel.Columns.Item(2).ColumnWidth = 22.1796875

# --- Sheet "Mes actual emp" (sheet2) ---
$wsEmp = $wb.Worksheets.Item("Mes actual emp")
$wsEmp.Range("B2").Value = "Christian Castaño Rodríguez"

# --- Sheet "Login" (sheet4) ---
$wsLogin = $wb.Worksheets.Item("Login")

# Update selections on each sheet (order chosen so the final active
# sheet/tab matches the target: "Mes actual emp" ends up selected).
[void]$wsCel.Activate()
[void]$wsCel.Range("D16").Select()

[void]$wsLogin.Activate()
[void]$wsLogin.Range("G16").Select()

[void]$wsEmp.Activate()
[void]$wsEmp.Range("F10").Select()
